$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comp_type_dmg_algo")

# ---------------------------------------------------------------------------
# The "damage_state_def" model table (comp_type_dmg_algo sheet) gained two
# new data columns: "location" (right after "median") and
# "recovery_function" (right before "recovery_mean"). Insert them as real
# column-inserts so everything to their right shifts over automatically.
# ---------------------------------------------------------------------------

# 1) Insert "location" column before column G (old G = "beta").
$ws.Columns.Item(7).Insert()

# 2) Insert "recovery_function" column before column N (old "recovery_mean",
#    already shifted right by the first insert).
$ws.Columns.Item(14).Insert()

# Header row (write recovery_function/Normal first so the shared-string
# table picks up the same ordering as the authoritative edit)
$ws.Range("N1").Value = "recovery_function"
$ws.Range("N2:N29").Value = "Normal"
$ws.Range("G1").Value = "location"

# Data rows 2-29: "location" defaults to 0, "recovery_function" to "Normal"
$ws.Range("G2:G29").Value = 0

# ---------------------------------------------------------------------------
# Column D ("is_piecewise") had a stray top border on every row; tidy it up
# so it matches the border pattern of its neighbouring C/E cells (keeps the
# boxed-group look only where C/E also have it).
# ---------------------------------------------------------------------------
$fixRows = 3,4,5,7,8,9,11,12,13,15,16,17,19,20,21,22,23,24,25,26,27,28,29
foreach ($r in $fixRows) {
    $ws.Range("D$r").Borders.Item(8).LineStyle = -4142
}

Write-Host "done"
